# Commit: "Sheet choice for #14"
# Adds a second worksheet ("Sheet2") positioned after the existing "1_data"
# sheet, populates it with a near-duplicate of the "1_data" data (the
# "skipped" shared-string label lands on G3 instead of E2/G2), makes it the
# active/selected sheet with G4 selected, and leaves "1_data" no longer the
# selected tab.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("1_data")

# Insert the new sheet right after "1_data" -- this also flips the
# workbook's activeTab to the new sheet's index automatically.
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# --- populate Sheet2 with the same headers/shared strings as "1_data" ---
$headers = @("id", "foo_a", "foo_b", "foo_c", "foo_d", "bar_a", "bar_b", "bar_c", "bar_d")
for ($col = 1; $col -le $headers.Length; $col++) {
    $ws2.Cells.Item(1, $col).Value = $headers[$col - 1]
}

$row2 = @(1, 6, 4, 4, 1, 4, 1, 3, 5)
for ($col = 1; $col -le $row2.Length; $col++) {
    $ws2.Cells.Item(2, $col).Value = $row2[$col - 1]
}

$row3 = @(2, 1, 4, 5, 3, 5, "skipped", 6, 1)
for ($col = 1; $col -le $row3.Length; $col++) {
    $ws2.Cells.Item(3, $col).Value = $row3[$col - 1]
}

$ws2.Cells.Item(4, 1).Value = 3
$ws2.Cells.Item(4, 8).Value = 4

# Make Sheet2 the active sheet/tab, with G4 as the selected cell.
[void]$ws2.Activate()
[void]$ws2.Range("G4").Select()
